# Auto-generated edit script: updates market-data columns (H-N) on each sheet
# to match the refreshed values from the scheduled Universalis price-sync run.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (132 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 5682727.5
$ws.Range("I12").Value = 22727772
$ws.Range("K12").Value = 22727772
$ws.Range("M12").Value = -22727602
$ws.Range("H13").Value = 6924.8335
$ws.Range("J13").Value = 6924.8335
$ws.Range("L13").Value = 6924.8335
$ws.Range("N13").Value = -7262.8335
$ws.Range("H17").Value = 160501.64
$ws.Range("J17").Value = 164391.92
$ws.Range("L17").Value = 493175.76
$ws.Range("N17").Value = -493511.76
$ws.Range("H40").Value = 1999.7693
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1999.7273
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1999.7273
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2349.7273
$ws.Range("H74").Value = 6616.8
$ws.Range("I74").Value = 6616.8
$ws.Range("K74").Value = 6616.8
$ws.Range("M74").Value = -5680.8
$ws.Range("H76").Value = 4749
$ws.Range("I76").Value = 4499
$ws.Range("J76").Value = 4999
$ws.Range("K76").Value = 4499
$ws.Range("L76").Value = 4999
$ws.Range("M76").Value = -4184
$ws.Range("N76").Value = -5629
$ws.Range("H77").Value = 6616.8
$ws.Range("I77").Value = 6616.8
$ws.Range("K77").Value = 33084
$ws.Range("M77").Value = -28404
$ws.Range("H79").Value = 4749
$ws.Range("I79").Value = 4499
$ws.Range("J79").Value = 4999
$ws.Range("K79").Value = 4499
$ws.Range("L79").Value = 4999
$ws.Range("M79").Value = -3407
$ws.Range("N79").Value = -7183
$ws.Range("H80").Value = 13158774
$ws.Range("I80").Value = 616.7
$ws.Range("J80").Value = 27778948
$ws.Range("K80").Value = 1850.1
$ws.Range("L80").Value = 83336844
$ws.Range("M80").Value = -852.1000000000001
$ws.Range("N80").Value = -83338840
$ws.Range("H83").Value = 13158774
$ws.Range("I83").Value = 616.7
$ws.Range("J83").Value = 27778948
$ws.Range("K83").Value = 5550.3
$ws.Range("L83").Value = 250010532
$ws.Range("M83").Value = -558.3000000000002
$ws.Range("N83").Value = -250020516
$ws.Range("H88").Value = 5182
$ws.Range("I88").Value = 1395.5
$ws.Range("J88").Value = 6444.1665
$ws.Range("K88").Value = 1395.5
$ws.Range("L88").Value = 6444.1665
$ws.Range("M88").Value = -989.5
$ws.Range("N88").Value = -7256.1665
$ws.Range("H91").Value = 5182
$ws.Range("I91").Value = 1395.5
$ws.Range("J91").Value = 6444.1665
$ws.Range("K91").Value = 1395.5
$ws.Range("L91").Value = 6444.1665
$ws.Range("M91").Value = 8.5
$ws.Range("N91").Value = -9252.166499999999
$ws.Range("H92").Value = 2038966.5
$ws.Range("I92").Value = 1116998.4
$ws.Range("J92").Value = 3473139
$ws.Range("K92").Value = 1116998.4
$ws.Range("L92").Value = 3473139
$ws.Range("M92").Value = -1115750.4
$ws.Range("N92").Value = -3475635
$ws.Range("H100").Value = 2433.7273
$ws.Range("I100").Value = 2433.7273
$ws.Range("K100").Value = 2433.7273
$ws.Range("M100").Value = -1892.7273
$ws.Range("H101").Value = 2192.8
$ws.Range("I101").Value = 2584
$ws.Range("J101").Value = 1932
$ws.Range("K101").Value = 7752
$ws.Range("L101").Value = 5796
$ws.Range("M101").Value = -6130
$ws.Range("N101").Value = -9040
$ws.Range("H107").Value = 16667020
$ws.Range("I107").Value = 393.4375
$ws.Range("K107").Value = 393.4375
$ws.Range("M107").Value = 1526.5625
$ws.Range("H112").Value = 3805.5
$ws.Range("I112").Value = 2749.1667
$ws.Range("J112").Value = 4157.6113
$ws.Range("K112").Value = 8247.500100000001
$ws.Range("L112").Value = 12472.8339
$ws.Range("M112").Value = -7139.500100000001
$ws.Range("N112").Value = -14688.8339
$ws.Range("H113").Value = 3568.8
$ws.Range("I113").Value = 2925
$ws.Range("K113").Value = 2925
$ws.Range("M113").Value = 329
$ws.Range("H116").Value = 12256.077
$ws.Range("J116").Value = 4337.5
$ws.Range("L116").Value = 4337.5
$ws.Range("N116").Value = -11221.5
$ws.Range("H121").Value = 1432.3334
$ws.Range("J121").Value = 1432.3334
$ws.Range("L121").Value = 4297.0002
$ws.Range("N121").Value = -7791.0002
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H137").Value = 4546710.5
$ws.Range("I137").Value = 1242.5333
$ws.Range("K137").Value = 3727.5999
$ws.Range("M137").Value = -1177.5999
$ws.Range("H138").Value = 5611.0566
$ws.Range("I138").Value = 14157.909
$ws.Range("J138").Value = 3372.5952
$ws.Range("K138").Value = 42473.727
$ws.Range("L138").Value = 10117.7856
$ws.Range("M138").Value = -37333.727
$ws.Range("N138").Value = -20397.7856
$ws.Range("H141").Value = 1414.4286
$ws.Range("I141").Value = 1138.6154
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 3415.8462
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 1764.1538
$ws.Range("N141").Value = -25360
$ws.Range("N123").ClearContents()

# --- Sheet: ARM (65 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 452.66666
$ws.Range("I5").Value = 488.30768
$ws.Range("K5").Value = 488.30768
$ws.Range("M5").Value = -376.30768
$ws.Range("H32").Value = 180062
$ws.Range("I32").Value = 193057.64
$ws.Range("K32").Value = 193057.64
$ws.Range("M32").Value = -192770.64
$ws.Range("H33").Value = 50005000
$ws.Range("I33").Value = 50005000
$ws.Range("K33").Value = 50005000
$ws.Range("M33").Value = -50004671
$ws.Range("H45").Value = 114076.22
$ws.Range("I45").Value = 127835.125
$ws.Range("K45").Value = 127835.125
$ws.Range("M45").Value = -127458.125
$ws.Range("H52").Value = 39999.4
$ws.Range("I52").Value = 39998.5
$ws.Range("K52").Value = 39998.5
$ws.Range("M52").Value = -39680.5
$ws.Range("H61").Value = 1270626.2
$ws.Range("I61").Value = 33657.273
$ws.Range("K61").Value = 33657.273
$ws.Range("M61").Value = -33445.273
$ws.Range("H74").Value = 344929.75
$ws.Range("I74").Value = 1489.619
$ws.Range("J74").Value = 972081.25
$ws.Range("K74").Value = 1489.619
$ws.Range("L74").Value = 972081.25
$ws.Range("M74").Value = -615.6189999999999
$ws.Range("N74").Value = -973829.25
$ws.Range("H77").Value = 344929.75
$ws.Range("I77").Value = 1489.619
$ws.Range("J77").Value = 972081.25
$ws.Range("K77").Value = 7448.094999999999
$ws.Range("L77").Value = 4860406.25
$ws.Range("M77").Value = -3080.094999999999
$ws.Range("N77").Value = -4869142.25
$ws.Range("H96").Value = 46883.168
$ws.Range("J96").Value = 46883.168
$ws.Range("L96").Value = 46883.168
$ws.Range("N96").Value = -52375.168
$ws.Range("H97").Value = 5658.4736
$ws.Range("I97").Value = 5950.6665
$ws.Range("K97").Value = 5950.6665
$ws.Range("M97").Value = -5454.6665
$ws.Range("H110").Value = 2980
$ws.Range("I110").Value = 2980
$ws.Range("K110").Value = 2980
$ws.Range("M110").Value = -935
$ws.Range("H122").Value = 1042.2
$ws.Range("I122").Value = 1047
$ws.Range("K122").Value = 3141
$ws.Range("M122").Value = -691
$ws.Range("H132").Value = 2986.6943
$ws.Range("I132").Value = 2563.889
$ws.Range("J132").Value = 4255.1113
$ws.Range("K132").Value = 7691.667
$ws.Range("L132").Value = 12765.3339
$ws.Range("M132").Value = -5161.667
$ws.Range("N132").Value = -17825.3339
$ws.Range("H136").Value = 1270626.2
$ws.Range("I136").Value = 33657.273
$ws.Range("K136").Value = 100971.819
$ws.Range("M136").Value = -98421.819

# --- Sheet: BSM (45 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 452.66666
$ws.Range("I4").Value = 488.30768
$ws.Range("K4").Value = 488.30768
$ws.Range("M4").Value = -373.30768
$ws.Range("H20").Value = 959.73334
$ws.Range("I20").Value = 1019
$ws.Range("J20").Value = 892
$ws.Range("K20").Value = 1019
$ws.Range("L20").Value = 892
$ws.Range("M20").Value = -772
$ws.Range("N20").Value = -1386
$ws.Range("H36").Value = 1389.2
$ws.Range("I36").Value = 1486.5
$ws.Range("J36").Value = 1000
$ws.Range("K36").Value = 1486.5
$ws.Range("L36").Value = 1000
$ws.Range("M36").Value = -952.5
$ws.Range("N36").Value = -2068
$ws.Range("H86").Value = 2239.0667
$ws.Range("I86").Value = 1217
$ws.Range("K86").Value = 1217
$ws.Range("M86").Value = -94
$ws.Range("H89").Value = 2239.0667
$ws.Range("I89").Value = 1217
$ws.Range("K89").Value = 6085
$ws.Range("M89").Value = -469
$ws.Range("H99").Value = 7912.9
$ws.Range("I99").Value = 10332.786
$ws.Range("J99").Value = 2266.5
$ws.Range("K99").Value = 10332.786
$ws.Range("L99").Value = 2266.5
$ws.Range("M99").Value = -8834.786
$ws.Range("N99").Value = -5262.5
$ws.Range("H101").Value = 8100
$ws.Range("I101").Value = 8100
$ws.Range("K101").Value = 8100
$ws.Range("M101").Value = -4855
$ws.Range("H132").Value = 103979.6
$ws.Range("H134").Value = 23079182
$ws.Range("I134").Value = 2005.6072
$ws.Range("J134").Value = 81821090
$ws.Range("K134").Value = 6016.821599999999
$ws.Range("L134").Value = 245463270
$ws.Range("M134").Value = -3481.821599999999
$ws.Range("N134").Value = -245468340

# --- Sheet: CRP (68 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 164.23077
$ws.Range("I7").Value = 206.1
$ws.Range("J7").Value = 24.666666
$ws.Range("K7").Value = 206.1
$ws.Range("L7").Value = 24.666666
$ws.Range("M7").Value = -93.09999999999999
$ws.Range("N7").Value = -250.666666
$ws.Range("H31").Value = 2512.1724
$ws.Range("I31").Value = 2000.6666
$ws.Range("J31").Value = 2957.6775
$ws.Range("K31").Value = 2000.6666
$ws.Range("L31").Value = 2957.6775
$ws.Range("M31").Value = -1705.6666
$ws.Range("N31").Value = -3547.6775
$ws.Range("H32").Value = 2299.6667
$ws.Range("I32").Value = 2299.6667
$ws.Range("K32").Value = 2299.6667
$ws.Range("M32").Value = -1983.6667
$ws.Range("H34").Value = 2512.1724
$ws.Range("I34").Value = 2000.6666
$ws.Range("J34").Value = 2957.6775
$ws.Range("K34").Value = 2000.6666
$ws.Range("L34").Value = 2957.6775
$ws.Range("M34").Value = -1798.6666
$ws.Range("N34").Value = -3361.6775
$ws.Range("H58").Value = 3822.889
$ws.Range("I58").Value = 3882
$ws.Range("K58").Value = 3882
$ws.Range("M58").Value = -3679
$ws.Range("H68").Value = 49995
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("H71").Value = 49995
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H107").Value = 1897.96
$ws.Range("I107").Value = 1563.8334
$ws.Range("K107").Value = 1563.8334
$ws.Range("M107").Value = 356.1666
$ws.Range("H132").Value = 25816.785
$ws.Range("I132").Value = 32085.727
$ws.Range("J132").Value = 2830.6667
$ws.Range("K132").Value = 96257.181
$ws.Range("L132").Value = 8492.000100000001
$ws.Range("M132").Value = -93727.181
$ws.Range("N132").Value = -13552.0001
$ws.Range("H134").Value = 2698.077
$ws.Range("I134").Value = 2094.5
$ws.Range("K134").Value = 6283.5
$ws.Range("M134").Value = -3748.5
$ws.Range("H136").Value = 3822.889
$ws.Range("I136").Value = 3882
$ws.Range("K136").Value = 11646
$ws.Range("M136").Value = -9096
$ws.Range("H141").Value = 332875.5
$ws.Range("J141").Value = 332875.5
$ws.Range("L141").Value = 332875.5
$ws.Range("N141").Value = -343235.5
$ws.Range("M68").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# --- Sheet: CUL (73 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 563.4
$ws.Range("I2").Value = 347.1
$ws.Range("K2").Value = 2082.6
$ws.Range("M2").Value = -1969.6
$ws.Range("H8").Value = 1350
$ws.Range("I8").Value = 1350
$ws.Range("K8").Value = 4050
$ws.Range("M8").Value = -3911
$ws.Range("H10").Value = 646.0714
$ws.Range("I10").Value = 9
$ws.Range("K10").Value = 27
$ws.Range("M10").Value = 112
$ws.Range("H12").Value = 34191.2
$ws.Range("I12").Value = 26.5
$ws.Range("J12").Value = 39447.31
$ws.Range("K12").Value = 79.5
$ws.Range("L12").Value = 118341.93
$ws.Range("M12").Value = 93.5
$ws.Range("N12").Value = -118687.93
$ws.Range("H55").Value = 2896.75
$ws.Range("I55").Value = 1308.8334
$ws.Range("J55").Value = 3849.5
$ws.Range("K55").Value = 3926.5002
$ws.Range("L55").Value = 11548.5
$ws.Range("M55").Value = -3749.5002
$ws.Range("N55").Value = -11902.5
$ws.Range("H58").Value = 11295.546
$ws.Range("I58").Value = 2648.3333
$ws.Range("J58").Value = 14538.25
$ws.Range("K58").Value = 7944.999899999999
$ws.Range("L58").Value = 43614.75
$ws.Range("M58").Value = -7816.999899999999
$ws.Range("N58").Value = -43870.75
$ws.Range("H68").Value = 2085.8
$ws.Range("I68").Value = 1711.6
$ws.Range("J68").Value = 2460
$ws.Range("K68").Value = 5134.799999999999
$ws.Range("L68").Value = 7380
$ws.Range("M68").Value = -4323.799999999999
$ws.Range("N68").Value = -9002
$ws.Range("H71").Value = 2085.8
$ws.Range("I71").Value = 1711.6
$ws.Range("J71").Value = 2460
$ws.Range("K71").Value = 15404.4
$ws.Range("L71").Value = 22140
$ws.Range("M71").Value = -11348.4
$ws.Range("N71").Value = -30252
$ws.Range("H105").Value = 18498
$ws.Range("J105").Value = 20622.5
$ws.Range("L105").Value = 61867.5
$ws.Range("N105").Value = -67109.5
$ws.Range("H113").Value = 390.33334
$ws.Range("I113").Value = 354
$ws.Range("J113").Value = 399.41666
$ws.Range("K113").Value = 1062
$ws.Range("L113").Value = 1198.24998
$ws.Range("M113").Value = 1108
$ws.Range("N113").Value = -5538.249980000001
$ws.Range("H121").Value = 6615.8335
$ws.Range("I121").Value = 564.6667
$ws.Range("K121").Value = 1694.0001
$ws.Range("M121").Value = -384.0001
$ws.Range("H132").Value = 1667.1666
$ws.Range("I132").Value = 1002
$ws.Range("J132").Value = 1999.75
$ws.Range("K132").Value = 9018
$ws.Range("L132").Value = 17997.75
$ws.Range("M132").Value = -6488
$ws.Range("N132").Value = -23057.75
$ws.Range("H140").Value = 3114.111
$ws.Range("I140").Value = 2253.5
$ws.Range("K140").Value = 6760.5
$ws.Range("M140").Value = -1580.5

# --- Sheet: GSM (40 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 64997
$ws.Range("J42").Value = 64997
$ws.Range("L42").Value = 64997
$ws.Range("N42").Value = -65967
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("H97").Value = 445.43243
$ws.Range("J97").Value = 560.0833
$ws.Range("L97").Value = 560.0833
$ws.Range("N97").Value = -1552.0833
$ws.Range("H104").Value = 52799.6
$ws.Range("J104").Value = 52799.6
$ws.Range("L104").Value = 52799.6
$ws.Range("N104").Value = -59787.6
$ws.Range("H106").Value = 39666.332
$ws.Range("J106").Value = 39666.332
$ws.Range("L106").Value = 39666.332
$ws.Range("N106").Value = -42190.332
$ws.Range("H107").Value = 53234.42
$ws.Range("I107").Value = 125245.5
$ws.Range("K107").Value = 125245.5
$ws.Range("M107").Value = -123325.5
$ws.Range("H115").Value = 64997
$ws.Range("J115").Value = 64997
$ws.Range("L115").Value = 64997
$ws.Range("N115").Value = -67347
$ws.Range("H122").Value = 1963.3871
$ws.Range("J122").Value = 2096.2222
$ws.Range("L122").Value = 6288.6666
$ws.Range("N122").Value = -11188.6666
$ws.Range("H132").Value = 422344.38
$ws.Range("I132").Value = 1480.6774
$ws.Range("K132").Value = 4442.0322
$ws.Range("M132").Value = -1912.0322
$ws.Range("H134").Value = 44142.855
$ws.Range("J134").Value = 44142.855
$ws.Range("L134").Value = 132428.565
$ws.Range("N134").Value = -137498.565
$ws.Range("N94").ClearContents()

# --- Sheet: LTW (71 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 9800
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H61").Value = 4568.2
$ws.Range("I61").Value = 3709
$ws.Range("J61").Value = 8005
$ws.Range("K61").Value = 3709
$ws.Range("L61").Value = 8005
$ws.Range("M61").Value = -3507
$ws.Range("N61").Value = -8409
$ws.Range("H82").Value = 793.3333
$ws.Range("I82").Value = 500.33334
$ws.Range("J82").Value = 891
$ws.Range("K82").Value = 500.33334
$ws.Range("L82").Value = 891
$ws.Range("M82").Value = -139.33334
$ws.Range("N82").Value = -1613
$ws.Range("H85").Value = 793.3333
$ws.Range("I85").Value = 500.33334
$ws.Range("J85").Value = 891
$ws.Range("K85").Value = 500.33334
$ws.Range("L85").Value = 891
$ws.Range("M85").Value = 747.66666
$ws.Range("N85").Value = -3387
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("H93").Value = 2377.2144
$ws.Range("I93").Value = 1273.5
$ws.Range("K93").Value = 1273.5
$ws.Range("M93").Value = -25.5
$ws.Range("H94").Value = 39325
$ws.Range("J94").Value = 39325
$ws.Range("L94").Value = 39325
$ws.Range("N94").Value = -40677
$ws.Range("H106").Value = 6946.3335
$ws.Range("J106").Value = 6946.3335
$ws.Range("L106").Value = 6946.3335
$ws.Range("N106").Value = -9470.333500000001
$ws.Range("H113").Value = 4568.2
$ws.Range("I113").Value = 3709
$ws.Range("J113").Value = 8005
$ws.Range("K113").Value = 3709
$ws.Range("L113").Value = 8005
$ws.Range("M113").Value = -1539
$ws.Range("N113").Value = -12345
$ws.Range("H122").Value = 2883.484
$ws.Range("I122").Value = 2774.8235
$ws.Range("J122").Value = 3015.4285
$ws.Range("K122").Value = 8324.470499999999
$ws.Range("L122").Value = 9046.2855
$ws.Range("M122").Value = -5874.470499999999
$ws.Range("N122").Value = -13946.2855
$ws.Range("H132").Value = 3290.5312
$ws.Range("I132").Value = 3040.8125
$ws.Range("K132").Value = 9122.4375
$ws.Range("M132").Value = -6592.4375
$ws.Range("H135").Value = 154996.67
$ws.Range("J135").Value = 154996.67
$ws.Range("L135").Value = 154996.67
$ws.Range("N135").Value = -165136.67
$ws.Range("H140").Value = 119989.5
$ws.Range("J140").Value = 119989.5
$ws.Range("L140").Value = 119989.5
$ws.Range("N140").Value = -130349.5
$ws.Range("N18").ClearContents()
$ws.Range("M88").ClearContents()
$ws.Range("M91").ClearContents()

# --- Sheet: WVR (31 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 95246.37
$ws.Range("J81").Value = 203159.4
$ws.Range("L81").Value = 406318.8
$ws.Range("N81").Value = -408440.8
$ws.Range("H84").Value = 95246.37
$ws.Range("J84").Value = 203159.4
$ws.Range("L84").Value = 2031594
$ws.Range("N84").Value = -2042202
$ws.Range("H107").Value = 866334.9
$ws.Range("I107").Value = 559.3333
$ws.Range("K107").Value = 1677.9999
$ws.Range("M107").Value = 242.0001
$ws.Range("H113").Value = 606.7222
$ws.Range("I113").Value = 527.6
$ws.Range("K113").Value = 1582.8
$ws.Range("M113").Value = 587.1999999999998
$ws.Range("H132").Value = 2099.5762
$ws.Range("I132").Value = 1698.0541
$ws.Range("K132").Value = 5094.1623
$ws.Range("M132").Value = -2564.1623
$ws.Range("H135").Value = 112530.664
$ws.Range("J135").Value = 112530.664
$ws.Range("L135").Value = 112530.664
$ws.Range("N135").Value = -122670.664
$ws.Range("H136").Value = 36864.715
$ws.Range("I136").Value = 50422.4
$ws.Range("J136").Value = 2970.5
$ws.Range("K136").Value = 151267.2
$ws.Range("L136").Value = 8911.5
$ws.Range("M136").Value = -148717.2
$ws.Range("N136").Value = -14011.5
